# Update cryptos list - Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.402.73"
$ws.Range("E2").Value = "  +6.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.349.59"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.45"
$ws.Range("E5").Value = "  +3.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.99"
$ws.Range("E6").Value = "  +2.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.586"
$ws.Range("E7").Value = "  +5.31%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.634"
$ws.Range("E9").Value = "  +1.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.16"
$ws.Range("E10").Value = "  +2.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0991"
$ws.Range("E11").Value = "  +3.54%  "

$ws.Range("E12").Value = "  +1.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.872.60"
$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.54"
$ws.Range("E14").Value = "  +5.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.39"
$ws.Range("E15").Value = "  +1.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.334.61"
$ws.Range("E16").Value = "  +2.37%  "

$ws.Range("E17").Value = "  +1.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.046.49"
$ws.Range("E18").Value = "  +5.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.74"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.36"
$ws.Range("E20").Value = "  +1.84%  "

$ws.Range("E21").Value = "  +5.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.13"
$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "305.49"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.83"
$ws.Range("E24").Value = "  +1.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.20"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.184"
$ws.Range("E26").Value = "  +9.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.68"
$ws.Range("E27").Value = "  +2.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.49"
$ws.Range("E28").Value = "  +2.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.90"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.48"
$ws.Range("E30").Value = "  +3.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").Value = "  +25.17%  "

$ws.Range("E32").Value = "  +4.63%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.58"
$ws.Range("E34").Value = "  +5.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.89"
$ws.Range("E35").Value = "  +6.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0508"
$ws.Range("E36").Value = "  +5.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.86"
$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("E38").Value = "  +1.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("E40").Value = "  -4.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.60"
$ws.Range("E41").Value = "  +3.38%  "

$ws.Range("E42").Value = "  +2.81%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.285"

$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.91"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.23"
$ws.Range("E48").Value = "  +1.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.206.60"
$ws.Range("E49").Value = "  +2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.09"
$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("E51").Value = "  +0.38%  "
